$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly data row needs to be inserted. The sheet already holds one
# row per (market, date) observation; a new observation is added right
# before the existing row 99, which pushes rows 99-126 down to 100-127
# (row 127 ends up holding what used to be row 126's data).
$ws.Rows("99").Insert()

# Populate the newly inserted row 99 with the new weekly observation.
$ws.Cells.Item(99, 1).Value  = 5
$ws.Cells.Item(99, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(99, 3).Value  = "Maule"
$ws.Cells.Item(99, 4).Value  = 44722
$ws.Cells.Item(99, 5).Value  = 7
$ws.Cells.Item(99, 6).Value  = 100112001
$ws.Cells.Item(99, 7).Value  = "Berenjena"
$ws.Cells.Item(99, 8).Value  = "Sin especificar"
$ws.Cells.Item(99, 9).Value  = "Primera"
$ws.Cells.Item(99, 10).Value = 300
$ws.Cells.Item(99, 11).Value = 6000
$ws.Cells.Item(99, 12).Value = 6000
$ws.Cells.Item(99, 13).Value = 6000
$ws.Cells.Item(99, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(99, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(99, 16).Value = 120
$ws.Cells.Item(99, 17).Value = 50
$ws.Cells.Item(99, 18).Value = "Hortaliza"

# Keep the date cell formatted like the rest of column D.
$ws.Cells.Item(99, 4).NumberFormat = $ws.Cells.Item(100, 4).NumberFormat
